# REG-05 "Verify required fields cannot be empty" test cases are added to the
# Registration Module sheet (rows 9-10), mirroring the existing REG-04 block
# (rows 7-8) in formatting: merged B/C columns, the same cell styles, and the
# same row-height pattern (31.5 / 47.25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration Module")

# --- Copy the formatting of the REG-04 block (rows 7-8) down onto the new
#     REG-05 block (rows 9-10) before writing any values, so the new rows end
#     up with the same cell styles (borders/fills/fonts/alignment) as the
#     existing rows. ---
$ws.Range("B7:L8").Copy()
$ws.Range("B9:L10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# --- Row 9: "Register with empty email" ---
$ws.Range("B9").Value = "REG-05"
$ws.Range("C9").Value = "Verify required fields cannot be empty"
$ws.Range("D9").Value = "REG-TC-06"
$ws.Range("E9").Value = "Registration Module"
$ws.Range("F9").Value = "REG-05"
$ws.Range("G9").Value = "Register with empty email"
$ws.Range("H9").Value = "On registration page"
$ws.Range("I9").Value = "Leave email empty, fill password, click register"
$ws.Range("J9").Value = '"" / Password123'
$ws.Range("K9").Value = "Error message appears"
$ws.Range("L9").Value = "High"

# --- Row 10: "Fill email, leave password empty, click register" ---
$ws.Range("D10").Value = "REG-TC-07"
$ws.Range("E10").Value = "Registration Module"
$ws.Range("F10").Value = "REG-05"
$ws.Range("G10").Value = "Fill email, leave password empty, click register"
$ws.Range("H10").Value = "On registration page"
$ws.Range("I10").Value = "Fill email, leave password empty, click registe"
$ws.Range("J10").Value = "user2@mail.com`n / " + [char]34 + [char]34
$ws.Range("K10").Value = "Error message appears"
$ws.Range("L10").Value = "High"

# --- Row heights to match the REG-04 block (wrapped-text auto height). ---
$ws.Rows.Item(9).RowHeight = 31.5
$ws.Rows.Item(10).RowHeight = 47.25

# --- Merge the Scenario ID / Scenario Description columns across the two new
#     rows, exactly like B7:B8 / C7:C8 above them. ---
$ws.Range("B9:B10").Merge()
$ws.Range("C9:C10").Merge()

# --- Hyperlinks on the "Test Data" column, matching the style of the rest of
#     the table: J9 links to the standard test mailbox (its display text
#     differs from the cell text, so "display" is written out); J10's link
#     target equals its own cell text, so no extra "display" is written. ---
$ws.Hyperlinks.Add($ws.Range("J9"), "mailto:user1@gmail.com", [Type]::Missing, [Type]::Missing, "user1@gmail.com")
$ws.Range("J9").Value = '"" / Password123'

$ws.Hyperlinks.Add($ws.Range("J10"), "mailto:user2@mail.com%0a / " + [char]34 + [char]34)
$ws.Range("J10").Value = "user2@mail.com`n / " + [char]34 + [char]34

# --- Move the selection/view to reflect where the author ended up editing. ---
$ws.Range("F5").Select()
$ws.Range("L10").Select()
